$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet to reflect the new extraction timestamp
$ws.Name = "IClientBalance-20241016-090632-"

# Bump the "Dt. Referencia" (column G) date serial for every data row
# (2024-10-15 -> 2024-10-16, i.e. serial 45580 -> 45581)
$ws.Range("G2:G274").Value = 45581

# Row 103: Saldo Previsto / Vl. Total corrected from 55771.97 to 999.99
$ws.Range("E103").Value = 999.99
$ws.Range("H103").Value = 999.99

# Row 104: Vl. Projetado set to 3303.53, Vl. Total recomputed to 8769.5
$ws.Range("D104").Value = 3303.53
$ws.Range("H104").Value = 8769.5

# Row 255: Saldo Previsto / Vl. Total corrected from 37351.61 to 999.99
$ws.Range("E255").Value = 999.99
$ws.Range("H255").Value = 999.99
